$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes: insert/delete rows to get from 42 to 45 rows ---
# 1) Insert a new row at the top -> row 1 becomes "Date and Time", everything else shifts down by one
$ws.Rows.Item(1).Insert()

# 2) The old "Maximum BMS Temperature in C" row (now row 32) is removed entirely
$ws.Rows.Item(32).Delete()

# 3) A new row is inserted (now row 35) for "Cycle Count of battery", before "Idling time percentage"
$ws.Rows.Item(35).Insert()

# --- Now (re)write every cell so the sheet matches the final layout exactly ---
$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-12 13:31:02.656000 to 2024-03-12 15:25:21.109000"
$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").Value = 0.07941324074074074
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"
$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = 35.43623027777778
$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = 1823.841569881111
$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = 4.135
$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = 4.423
$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 10
$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 99
$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = 56.7960429328498
$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = 32.11212393858929
$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = 89
$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Eco mode`n97.11%`nCustom mode`n1.97%`nSports mode`n0.83%"
$ws.Rows.Item(12).AutoFit()
$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = 6129.537488
$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = -961.6036396561219
$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = 134.2332089891667
$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 6.855366834693272
$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.376
$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 3.106
$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = 0.27
$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = 32
$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = 42
$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 10
$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = 58
$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = 52
$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = 51
$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = 48
$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = 0
$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = 0
$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 42
$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 32
$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = 10
$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 55
$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.948992665277778
$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.00000007890786349891408
$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 115
$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = 5.18963018963019
$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 10.28948528948529
$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = 7.397782397782398
$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = 11.45498645498645
$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = 44.86549486549487
$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = 18.42909342909343
$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = 0.2488502488502489
$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = 0.3858753858753859
$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 0.3087003087003087
$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0
